$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A51").Value = "red pepper"
$ws.Range("B51").Value = "Vegetable"
$ws.Range("C51:G51").Value = 0
$ws.Range("A52").Value = "green pepper"
$ws.Range("B52").Value = "Vegetable"
$ws.Range("C52:G52").Value = 0
$ws.Range("A53").Value = "aubergine"
$ws.Range("B53").Value = "Vegetable"
$ws.Range("C53:G53").Value = 0
$ws.Range("A54").Value = "courgette"
$ws.Range("B54").Value = "Vegetable"
$ws.Range("C54:G54").Value = 0
$ws.Range("A55").Value = "sunflower seeds"
$ws.Range("B55").Value = "Grocery"
$ws.Range("C55:G55").Value = 0
$ws.Range("A56").Value = "cottage cheese"
$ws.Range("B56").Value = "Dairy"
$ws.Range("C56:G56").Value = 0
$ws.Range("A57").Value = "sour cream"
$ws.Range("B57").Value = "Dairy"
$ws.Range("C57:G57").Value = 0
$ws.Range("A58").Value = "bacon"
$ws.Range("B58").Value = "Meat"
$ws.Range("C58:G58").Value = 0
$ws.Range("A59").Value = "chorizo sausage"
$ws.Range("B59").Value = "Meat"
$ws.Range("C59:G59").Value = 0
$ws.Range("A60").Value = "red lentils"
$ws.Range("B60").Value = "Grocery"
$ws.Range("C60:G60").Value = 0
$ws.Range("A61").Value = "mushrooms"
$ws.Range("B61").Value = "Vegetable"
$ws.Range("C61:G61").Value = 0
$ws.Range("A62").Value = "soy sauce"
$ws.Range("B62").Value = "Check"
$ws.Range("C62:G62").Value = 0
$ws.Range("A63").Value = "sherry"
$ws.Range("B63").Value = "Check"
$ws.Range("C63:G63").Value = 0
$ws.Range("A64").Value = "pearl barley"
$ws.Range("B64").Value = "Grocery"
$ws.Range("C64:G64").Value = 0
$ws.Range("A65").Value = "canned kidney beans"
$ws.Range("B65").Value = "Grocery"
$ws.Range("C65:G65").Value = 0

$ws.Range("A21").Value = "green chilli"
$ws.Range("A39").Value = "red chilli"

$ws.AutoFilterMode = $false
$ws.Range("A1:G65").AutoFilter()

foreach ($n in $wb.Names) {
    if ($n.Name -eq "tblIngredients!_FilterDatabase") {
        $n.RefersTo = "=tblIngredients!`$A`$1:`$G`$65"
    }
}

$ws.Range("A66").Value = "cumin seeds"
$ws.Range("B66").Value = "Check"
$ws.Range("C66:G66").Value = 0
$ws.Range("A67").Value = "miso"
$ws.Range("B67").Value = "Check"
$ws.Range("C67:G67").Value = 0
$ws.Range("A68").Value = "asafoetida"
$ws.Range("B68").Value = "Check"
$ws.Range("C68:G68").Value = 0

$ws.Range("A68").Select()
